$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 62

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 4).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2023-06-20"
$ws.Cells.Item($row, 2).Value = "16:16:53"
$ws.Cells.Item($row, 3).Value = "Tuesday"
$ws.Cells.Item($row, 4).Value = "25"
$ws.Cells.Item($row, 5).Value = 122141
$ws.Cells.Item($row, 6).Value = 133706
$ws.Cells.Item($row, 7).Value = 162431
$ws.Cells.Item($row, 8).Value = 133290
$ws.Cells.Item($row, 9).Value = 177326
$ws.Cells.Item($row, 10).Value = 114523
$ws.Cells.Item($row, 11).Value = 201562
$ws.Cells.Item($row, 12).Value = 225237
$ws.Cells.Item($row, 13).Value = 175501
$ws.Cells.Item($row, 14).Value = 103845
$ws.Cells.Item($row, 15).Value = 39217
$ws.Cells.Item($row, 16).Value = 33873
$ws.Cells.Item($row, 17).Value = 51902
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36176
$ws.Cells.Item($row, 20).Value = -1
